# Apply cryptos list price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers to Excel's type
# parser ("588.64", "0.534", etc.). Force them to Text format *before* writing
# so Excel stores the exact original string instead of a floating point
# approximation (and so short trailing digits like "0.740" aren't trimmed).
$textCells = @(
    "D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D29", "D31", "D33", "D34", "D36", "D38", "D39", "D40", "D41", "D43", "D46", "D47", "D48", "D49", "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row-by-row cell updates (coin name / link / price / 1h volume%)

# Row 2
$ws.Range("D2").Value = "62.993.86"
$ws.Range("E2").Value = "  +6.79%  "

# Row 3
$ws.Range("D3").Value = "3.121.03"
$ws.Range("E3").Value = "  +4.42%  "

# Row 4
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").Value = "588.64"
$ws.Range("E5").Value = "  +5.21%  "

# Row 6
$ws.Range("D6").Value = "144.65"
$ws.Range("E6").Value = "  +5.93%  "

# Row 7
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("D8").Value = "3.107.68"
$ws.Range("E8").Value = "  +4.31%  "

# Row 9
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  +2.36%  "

# Row 10
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  +10.16%  "

# Row 11
$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  +12.95%  "

# Row 12
$ws.Range("D12").Value = "0.472"
$ws.Range("E12").Value = "  +3.82%  "

# Row 13
$ws.Range("D13").Value = "0.0000245"
$ws.Range("E13").Value = "  +7.02%  "

# Row 14
$ws.Range("D14").Value = "35.78"
$ws.Range("E14").Value = "  +6.99%  "

# Row 15
$ws.Range("E15").Value = "  +1.08%  "

# Row 16
$ws.Range("D16").Value = "3.625.87"
$ws.Range("E16").Value = "  +4.03%  "

# Row 17
$ws.Range("D17").Value = "7.36"
$ws.Range("E17").Value = "  +1.05%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.107.00"
$ws.Range("E18").Value = "  +3.83%  "

# Row 19
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "62.864.66"
$ws.Range("E19").Value = "  +6.31%  "

# Row 20
$ws.Range("D20").Value = "455.25"
$ws.Range("E20").Value = "  +6.52%  "

# Row 21
$ws.Range("D21").Value = "14.18"
$ws.Range("E21").Value = "  +4.20%  "

# Row 22
$ws.Range("D22").Value = "0.740"
$ws.Range("E22").Value = "  +3.05%  "

# Row 23
$ws.Range("D23").Value = "7.54"
$ws.Range("E23").Value = "  +6.54%  "

# Row 24
$ws.Range("D24").Value = "13.82"
$ws.Range("E24").Value = "  +4.41%  "

# Row 25
$ws.Range("D25").Value = "82.54"
$ws.Range("E25").Value = "  +2.70%  "

# Row 27
$ws.Range("E27").Value = "  +5.21%  "

# Row 28
$ws.Range("D28").Value = "2.71"
$ws.Range("E28").Value = "  +7.02%  "

# Row 29
$ws.Range("D29").Value = "8.31"
$ws.Range("E29").Value = "  +7.24%  "

# Row 31
$ws.Range("D31").Value = "6.89"
$ws.Range("E31").Value = "  +14.98%  "

# Row 32
$ws.Range("E32").Value = "  +16.49%  "

# Row 33
$ws.Range("D33").Value = "27.27"
$ws.Range("E33").Value = "  +6.58%  "

# Row 34
$ws.Range("D34").Value = "1.05"
$ws.Range("E34").Value = "  +5.35%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0805"
$ws.Range("E35").Value = "  +7.39%  "

# Row 36
$ws.Range("D36").Value = "6.12"
$ws.Range("E36").Value = "  +3.78%  "

# Row 37
$ws.Range("E37").Value = "  +7.25%  "

# Row 38
$ws.Range("D38").Value = "50.66"
$ws.Range("E38").Value = "  +3.68%  "

# Row 39
$ws.Range("D39").Value = "3.04"
$ws.Range("E39").Value = "  +12.26%  "

# Row 40
$ws.Range("D40").Value = "8.89"
$ws.Range("E40").Value = "  +2.70%  "

# Row 41
$ws.Range("D41").Value = "427.78"
$ws.Range("E41").Value = "  +7.93%  "

# Row 42
$ws.Range("D42").Value = "2.939.22"
$ws.Range("E42").Value = "  +7.09%  "

# Row 43
$ws.Range("D43").Value = "0.0372"
$ws.Range("E43").Value = "  +6.26%  "

# Row 44
$ws.Range("E44").Value = "  +11.61%  "

# Row 45
$ws.Range("E45").Value = "  +2.55%  "

# Row 46
$ws.Range("D46").Value = "2.19"
$ws.Range("E46").Value = "  +10.35%  "

# Row 47
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "35.23"
$ws.Range("E47").Value = "  +0.44%  "

# Row 48
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.02%  "

# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "124.10"
$ws.Range("E49").Value = "  +1.26%  "

# Row 50
$ws.Range("E50").Value = "  +1.84%  "

# Row 51
$ws.Range("D51").Value = "24.79"
$ws.Range("E51").Value = "  +6.73%  "
